# Edit script for LOT2055.docx
#
# The body text of several paragraphs was permuted (moved around between
# paragraph "slots") while paragraph styles/formatting stayed put. This is
# implemented as two cycles of text moves. To do this safely with Word's
# Find & Replace (so that a later search never accidentally matches text
# that an earlier replacement just inserted), we first move every "before"
# text into a unique placeholder token, and only afterwards replace the
# placeholders with their final ("after") text.

$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $oldText)
    }
}

# ---- Phase 1: move every original text to a unique placeholder ----

Replace-Text "Desenvolver nos alunos as competências necessárias para aplicar conceitos e ferramentas de estatística em bioprocessos, com ênfase em planejamento e otimização de experimentos." "@@PLACEHOLDER_1@@"

Replace-Text "Develop in students the necessary skills to apply concepts and tools of statistics in bioprocesses, with an emphasis on experimental design and optimization." "@@PLACEHOLDER_2@@"

Replace-Text "5817181 - Valdeir Arantes" "@@PLACEHOLDER_3@@"

Replace-Text "1. Fundamentos de estatística aplicada; 2. Análise de sistemas de medição; 3. Análise de Variância; 4. Testes de comparações múltiplas; 5. Controle estatístico de processos; 6. Planejamento de Experimentos: planejamentos fatoriais, superfícies de resposta, planejamentos de mistura; 7. Aplicação de software estatístico e estratégia sequencial de planejamentos experimentais." "@@PLACEHOLDER_4@@"

Replace-Text "1. The role of statistics in Engineering; Fundamentals of applied statistics; Analysis of Variance; Multiple comparison tests; Experimental Design" "@@PLACEHOLDER_5@@"

Replace-Text "1. O papel da estatística na Engenharia: métodos de coleta de dados 2. Fundamentos de estatística aplicada 3. Análise de Variância: análise de variância de um modelo 4. Testes de comparações múltiplas (Tukey, Hsu) 5. Planejamento de Experimentos: vantagens dos experimentos fatoriais em relação aos experimentos do tipo um fator por vez; varielaboração do planejamento fatorial Completo do tipo 2^k e fracionado, e superfície de resposta" "@@PLACEHOLDER_6@@"

Replace-Text "A avaliação será composta por provas, exercícios, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n." "@@PLACEHOLDER_7@@"

Replace-Text "MF≥ 5,0 para aprovação 5,0" "@@PLACEHOLDER_8@@"

Replace-Text "(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada." "@@PLACEHOLDER_9@@"

Replace-Text "1. BOX, G.E.P.; HUNTER, W.G.; HUNTER, J.S. Statistics for Experimenters: an introduction to designs, data analysis and model building. New York: John Wiley & Sons Inc., 1978.^l^l2. RODRIGUES, M. I. e IEMMA, A. F. Planejamento de experimentos e otimização de processos. Campinas: Cárita editora, 2009.^l^l3. Planejamento e otimização de Experimentos. Roy E. Bruns, Edit. UNICAMP, 1996" "@@PLACEHOLDER_10@@"

# ---- Phase 2: move placeholders to their final ("after") text ----
# (placeholder N) -> (original text that used to belong to placeholder M)

Replace-Text "@@PLACEHOLDER_1@@" "1. Fundamentos de estatística aplicada; 2. Análise de sistemas de medição; 3. Análise de Variância; 4. Testes de comparações múltiplas; 5. Controle estatístico de processos; 6. Planejamento de Experimentos: planejamentos fatoriais, superfícies de resposta, planejamentos de mistura; 7. Aplicação de software estatístico e estratégia sequencial de planejamentos experimentais."

Replace-Text "@@PLACEHOLDER_2@@" "1. The role of statistics in Engineering; Fundamentals of applied statistics; Analysis of Variance; Multiple comparison tests; Experimental Design"

Replace-Text "@@PLACEHOLDER_3@@" "Desenvolver nos alunos as competências necessárias para aplicar conceitos e ferramentas de estatística em bioprocessos, com ênfase em planejamento e otimização de experimentos."

Replace-Text "@@PLACEHOLDER_4@@" "1. O papel da estatística na Engenharia: métodos de coleta de dados 2. Fundamentos de estatística aplicada 3. Análise de Variância: análise de variância de um modelo 4. Testes de comparações múltiplas (Tukey, Hsu) 5. Planejamento de Experimentos: vantagens dos experimentos fatoriais em relação aos experimentos do tipo um fator por vez; varielaboração do planejamento fatorial Completo do tipo 2^k e fracionado, e superfície de resposta"

Replace-Text "@@PLACEHOLDER_5@@" "Develop in students the necessary skills to apply concepts and tools of statistics in bioprocesses, with an emphasis on experimental design and optimization."

Replace-Text "@@PLACEHOLDER_6@@" "A avaliação será composta por provas, exercícios, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

Replace-Text "@@PLACEHOLDER_7@@" "MF≥ 5,0 para aprovação 5,0"

Replace-Text "@@PLACEHOLDER_8@@" "(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada."

Replace-Text "@@PLACEHOLDER_9@@" "1. BOX, G.E.P.; HUNTER, W.G.; HUNTER, J.S. Statistics for Experimenters: an introduction to designs, data analysis and model building. New York: John Wiley & Sons Inc., 1978.^l^l2. RODRIGUES, M. I. e IEMMA, A. F. Planejamento de experimentos e otimização de processos. Campinas: Cárita editora, 2009.^l^l3. Planejamento e otimização de Experimentos. Roy E. Bruns, Edit. UNICAMP, 1996"

Replace-Text "@@PLACEHOLDER_10@@" "5817181 - Valdeir Arantes"

Write-Output "DONE"
